# Refresh cryptos.xlsx price / volume snapshot (GitHub Actions scheduled job).
# Updates the Price (column D) and Volume(1h) (column E) cells for the
# data rows (2-51), and swaps the WOONetwork / Elrond rows (50-51), which
# changed rank order in this run. All Price/Volume cells are kept as plain
# text (matching the source sheet's inline-string cells), so NumberFormat
# is forced to Text before the assignment and the style is reset back to
# Normal afterwards to avoid leaving a visible formatting change behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# row -> new Price text (column D); only rows whose price actually changed.
$priceUpdates = @{
    2 = "30.580.08"
    3 = "1.915.60"
    4 = "1.0000"
    5 = "244.96"
    6 = "1.0000"
    7 = "0.4838"
    8 = "0.2894"
    9 = "0.06796"
    10 = "111.88"
    11 = "19.44"
    12 = "1.916.42"
    13 = "0.07572"
    14 = "5.403"
    15 = "0.6701"
    16 = "294.25"
    17 = "30.572.48"
    18 = "13.03"
    19 = "0.000007600"
    21 = "5.511"
    22 = "2.164.10"
    23 = "0.9998"
    24 = "6.417"
    25 = "9.461"
    26 = "166.10"
    28 = "2.076"
    29 = "0.1065"
    31 = "4.127"
    32 = "4.042"
    33 = "0.04983"
    34 = "0.7336"
    35 = "1.141"
    36 = "2.717"
    37 = "0.02024"
    39 = "2.018"
    40 = "109.26"
    41 = "0.4427"
    42 = "0.8657"
    43 = "5.840"
    45 = "69.33"
    46 = "7.207"
    47 = "48.66"
    48 = "9.228"
    49 = "0.1226"
}

# row -> new Volume(1h) text (column E); every data row 2-49 changed.
$volumeUpdates = @{
    2 = "  -0.10%  "
    3 = "  -0.43%  "
    4 = "  -0.05%  "
    5 = "  -0.74%  "
    6 = "  -0.05%  "
    7 = "  +1.88%  "
    8 = "  -0.05%  "
    9 = "  -0.72%  "
    10 = "  +6.26%  "
    11 = "  +5.64%  "
    12 = "  -0.42%  "
    13 = "  -1.50%  "
    14 = "  +1.25%  "
    15 = "  +0.23%  "
    16 = "  +1.75%  "
    17 = "  -0.18%  "
    18 = "  +0.61%  "
    19 = "  -0.19%  "
    21 = "  -1.21%  "
    22 = "  -0.35%  "
    23 = "  -0.14%  "
    24 = "  -0.68%  "
    25 = "  -0.59%  "
    26 = "  -0.55%  "
    27 = "  -4.28%  "
    28 = "  -1.85%  "
    29 = "  -0.62%  "
    30 = "  +2.80%  "
    31 = "  -1.09%  "
    32 = "  +0.02%  "
    33 = "  -0.88%  "
    34 = "  +0.37%  "
    35 = "  -0.15%  "
    36 = "  -0.66%  "
    37 = "  -2.00%  "
    38 = "  -0.24%  "
    39 = "  -1.51%  "
    40 = "  -2.23%  "
    41 = "  +0.77%  "
    42 = "  -0.74%  "
    43 = "  -1.56%  "
    44 = "  -0.03%  "
    45 = "  +2.38%  "
    46 = "  -1.16%  "
    47 = "  -0.41%  "
    48 = "  -1.16%  "
    49 = "  -1.44%  "
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $row 4 $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

# Rows 50/51 swapped rank order (WOONetwork <-> Elrond) with refreshed data.
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue 50 4 "34.75"
$ws.Cells.Item(50, 5).Value = "  -0.56%  "

$ws.Cells.Item(51, 2).Value = "WOONetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue 51 4 "0.2493"
$ws.Cells.Item(51, 5).Value = "  -0.39%  "
